# Update cryptocurrency price/volume figures on Sheet1 to reflect the
# latest GitHub Actions data refresh (commit: "Updated symbol list on
# Mon Dec 26 09:35:16 UTC 2022 with GitHub Actions").
#
# All of these cells are stored as text (not numbers) in the workbook,
# so we force the target cells to a text number format before writing
# the new values. This prevents Excel from reinterpreting numeric-
# looking strings (e.g. "243.22") as actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and E ("Volume(1h)") updates, keyed by cell address.
$updates = @{
    "D2"  = "243.22"
    "D3"  = "23.01"
    "D4"  = "5.405"
    "D5"  = "0.05920"
    "D7"  = "6.552"
    "D8"  = "0.8114"
    "D9"  = "0.9093"
    "D10" = "0.1406"
    "D11" = "0.07346"
    "D13" = "0.03046"
    "D14" = "0.09350"
    "D15" = "3.846"
    "D16" = "0.001562"
    "D17" = "0.04666"
    "D18" = "0.0005941"
    "D19" = "0.006082"
    "E20" = "19HotbitTokenHTB"
    "D21" = "0.0009810"
    "D22" = "0.00009402"
    "D23" = "3.606"
    "D24" = "2.138"
    "D27" = "0.0002900"
    "D40" = "0.03969"
    "D41" = "0.006197"
    "E41" = "40KickTokenKICKBestin24h"
    "D42" = "0.1073"
    "D43" = "0.003000"
    "D45" = "0.00005259"
    "D47" = "0.7821"
    "D48" = "0.002264"
    "D49" = "0.00002100"
    "D50" = "0.0002000"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
